# quarto 1.8 updates and fixes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the "example" link for Week 6 (row 7): bump example number 05 -> 06
$ws.Range("I7").Value = "/example/06-example-dplyr-and-more"

# Fix the "assignment" link for Week 7 (row 8): add missing leading slash
$ws.Range("J8").Value = "/assignment/07-assignment-fars"

# Update the active cell / selection to reflect where editing left off
$ws.Range("J9").Select()
